$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.088703393936157
$ws.Range("B1").Value = 2.514638900756836
$ws.Range("C1").Value = 2.621349573135376
$ws.Range("D1").Value = 3.242945194244385
$ws.Range("E1").Value = 0.8074414134025574
